$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new "MAE" column before the existing "Tipo" column (D)
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "MAE"
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Data row
$ws.Range("D2").Value = 0.1361288253571666

# Adjust the MSE value slightly (recalculated value from diff)
$ws.Range("B2").Value = 0.04215534119371403
